# Auto-generated: update Leve profit-tracking numbers (currentAveragePrice /
# LevePrice / LeveProfit columns H-N) across all eight crafting-job sheets,
# reflecting refreshed market-board data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4712.75
$ws.Range("I2").Value = 265
$ws.Range("K2").Value = 265
$ws.Range("M2").Value = -152
$ws.Range("H17").Value = 2082.72
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2082.72
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6248.16
$ws.Range("N17").Value = -6584.16
$ws.Range("M17").ClearContents()
$ws.Range("H112").Value = 3392
$ws.Range("J112").Value = 3723.4167
$ws.Range("L112").Value = 11170.2501
$ws.Range("N112").Value = -13386.2501
$ws.Range("H125").Value = 210254.53
$ws.Range("I125").Value = 11083.333
$ws.Range("J125").Value = 260047.33
$ws.Range("K125").Value = 99749.997
$ws.Range("L125").Value = 2340425.97
$ws.Range("M125").Value = -97289.997
$ws.Range("N125").Value = -2345345.97
$ws.Range("H137").Value = 4172.849
$ws.Range("I137").Value = 4021.2646
$ws.Range("J137").Value = 4305
$ws.Range("K137").Value = 12063.7938
$ws.Range("L137").Value = 12915
$ws.Range("M137").Value = -9513.793799999999
$ws.Range("N137").Value = -18015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1055.0278
$ws.Range("I2").Value = 1073.2069
$ws.Range("J2").Value = 979.7143
$ws.Range("K2").Value = 1073.2069
$ws.Range("L2").Value = 979.7143
$ws.Range("M2").Value = -960.2068999999999
$ws.Range("N2").Value = -1205.7143
$ws.Range("H32").Value = 31997
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 31997
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 31997
$ws.Range("N32").Value = -32571
$ws.Range("M32").ClearContents()
$ws.Range("H74").Value = 3916.7932
$ws.Range("I74").Value = 2083.3157
$ws.Range("K74").Value = 2083.3157
$ws.Range("M74").Value = -1209.3157
$ws.Range("H77").Value = 3916.7932
$ws.Range("I77").Value = 2083.3157
$ws.Range("K77").Value = 10416.5785
$ws.Range("M77").Value = -6048.5785
$ws.Range("H102").Value = 995
$ws.Range("I102").Value = 995
$ws.Range("K102").Value = 995
$ws.Range("M102").Value = 627
$ws.Range("H110").Value = 2682.2
$ws.Range("I110").Value = 2499.3333
$ws.Range("J110").Value = 2956.5
$ws.Range("K110").Value = 2499.3333
$ws.Range("L110").Value = 2956.5
$ws.Range("M110").Value = -454.3332999999998
$ws.Range("N110").Value = -7046.5
$ws.Range("H116").Value = 1055.0278
$ws.Range("I116").Value = 1073.2069
$ws.Range("J116").Value = 979.7143
$ws.Range("K116").Value = 1073.2069
$ws.Range("L116").Value = 979.7143
$ws.Range("M116").Value = 1220.7931
$ws.Range("N116").Value = -5567.7143
$ws.Range("H122").Value = 3726.3333
$ws.Range("I122").Value = 3726.3333
$ws.Range("K122").Value = 11178.9999
$ws.Range("M122").Value = -8728.999899999999
$ws.Range("H132").Value = 4136.8486
$ws.Range("I132").Value = 2810.4138
$ws.Range("J132").Value = 13753.5
$ws.Range("K132").Value = 8431.241399999999
$ws.Range("L132").Value = 41260.5
$ws.Range("M132").Value = -5901.241399999999
$ws.Range("N132").Value = -46320.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1055.0278
$ws.Range("I3").Value = 1073.2069
$ws.Range("J3").Value = 979.7143
$ws.Range("K3").Value = 1073.2069
$ws.Range("L3").Value = 979.7143
$ws.Range("M3").Value = -959.2068999999999
$ws.Range("N3").Value = -1207.7143
$ws.Range("H86").Value = 6387.357
$ws.Range("I86").Value = 3585.2856
$ws.Range("K86").Value = 3585.2856
$ws.Range("M86").Value = -2462.2856
$ws.Range("H89").Value = 6387.357
$ws.Range("I89").Value = 3585.2856
$ws.Range("K89").Value = 17926.428
$ws.Range("M89").Value = -12310.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2858.889
$ws.Range("I16").Value = 2788.5
$ws.Range("K16").Value = 2788.5
$ws.Range("M16").Value = -2501.5
$ws.Range("H31").Value = 5868.9546
$ws.Range("I31").Value = 2920.8
$ws.Range("K31").Value = 2920.8
$ws.Range("M31").Value = -2625.8
$ws.Range("H34").Value = 5868.9546
$ws.Range("I34").Value = 2920.8
$ws.Range("K34").Value = 2920.8
$ws.Range("M34").Value = -2718.8
$ws.Range("H113").Value = 2858.889
$ws.Range("I113").Value = 2788.5
$ws.Range("K113").Value = 2788.5
$ws.Range("M113").Value = -618.5
$ws.Range("H132").Value = 3602.721
$ws.Range("I132").Value = 3138.8
$ws.Range("J132").Value = 5632.375
$ws.Range("K132").Value = 9416.400000000001
$ws.Range("L132").Value = 16897.125
$ws.Range("M132").Value = -6886.400000000001
$ws.Range("N132").Value = -21957.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 856.9231
$ws.Range("I12").Value = 624
$ws.Range("J12").Value = 960.44446
$ws.Range("K12").Value = 1872
$ws.Range("L12").Value = 2881.33338
$ws.Range("M12").Value = -1699
$ws.Range("N12").Value = -3227.33338
$ws.Range("H74").Value = 6265.6665
$ws.Range("J74").Value = 10800
$ws.Range("L74").Value = 32400
$ws.Range("N74").Value = -34522
$ws.Range("H77").Value = 6265.6665
$ws.Range("J77").Value = 10800
$ws.Range("L77").Value = 97200
$ws.Range("N77").Value = -107808
$ws.Range("H97").Value = 3865.6667
$ws.Range("I97").Value = 3298.5
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 9895.5
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = -9399.5
$ws.Range("N97").Value = -15992
$ws.Range("H129").Value = 9935439
$ws.Range("J129").Value = 20973646
$ws.Range("L129").Value = 62920938
$ws.Range("N129").Value = -62930938

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2897.0908
$ws.Range("J102").Value = 7119.6
$ws.Range("L102").Value = 7119.6
$ws.Range("N102").Value = -10363.6
$ws.Range("H107").Value = 573.125
$ws.Range("I107").Value = 377
$ws.Range("K107").Value = 377
$ws.Range("M107").Value = 1543
$ws.Range("H113").Value = 6425.5454
$ws.Range("I113").Value = 5898.6665
$ws.Range("J113").Value = 7057.8
$ws.Range("K113").Value = 5898.6665
$ws.Range("L113").Value = 7057.8
$ws.Range("M113").Value = -3728.6665
$ws.Range("N113").Value = -11397.8
$ws.Range("H132").Value = 4539.3877
$ws.Range("J132").Value = 7159.9
$ws.Range("L132").Value = 21479.7
$ws.Range("N132").Value = -26539.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1563.0952
$ws.Range("I61").Value = 1563.0952
$ws.Range("K61").Value = 1563.0952
$ws.Range("M61").Value = -1361.0952
$ws.Range("H113").Value = 1563.0952
$ws.Range("I113").Value = 1563.0952
$ws.Range("K113").Value = 1563.0952
$ws.Range("M113").Value = 606.9048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 985.3611
$ws.Range("I113").Value = 357.92593
$ws.Range("J113").Value = 2867.6667
$ws.Range("K113").Value = 1073.77779
$ws.Range("L113").Value = 8603.000100000001
$ws.Range("M113").Value = 1096.22221
$ws.Range("N113").Value = -12943.0001
$ws.Range("H126").Value = 2087.3333
$ws.Range("I126").Value = 1797
$ws.Range("J126").Value = 3539
$ws.Range("K126").Value = 5391
$ws.Range("L126").Value = 10617
$ws.Range("M126").Value = -2921
$ws.Range("N126").Value = -15557
$ws.Range("H136").Value = 7822.067
$ws.Range("J136").Value = 6831.125
$ws.Range("L136").Value = 20493.375
$ws.Range("N136").Value = -25593.375
